$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63

# Column A's new value ("2025-04-29") is date-shaped text. A plain
# .Value assignment gets auto-converted into a date serial number (the
# same way real Excel parses typed-in date-looking strings), which would
# also tack on a date NumberFormat/style - neither of which the target
# row should have (every other row in column A is stored as literal
# text with the default style). Route the text through a formula that
# evaluates to the literal string, then flatten that formula down to a
# static value via copy / paste-values so it lands as plain text with no
# extra formatting.
$cellA = $ws.Range("A" + $row)
$cellA.Formula = '="2025-04-29"'
$cellA.Copy()
$cellA.PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B" + $row).Value = "ruissellement"
$ws.Range("C" + $row).Value = 39
$ws.Range("D" + $row).Value = 1
